$d = $word.ActiveDocument

# The document ends with two empty paragraphs (one right after the
# "Dependencies" text, and the final paragraph before the sectPr).
# We need to insert a new "Biography" paragraph and a paragraph
# containing a hyperlink to a bio-examples article, in between those
# two existing empty paragraphs - leaving both of them untouched.

$lastPara = $d.Paragraphs.Last
$anchor = $lastPara.Previous()

# Insert a new empty paragraph right after the anchor (the first of the
# two trailing empty paragraphs) and fill it with "Biography".
$anchor.Range.InsertParagraphAfter()
$bioPara = $anchor.Next()
$bioPara.Range.Text = "Biography"

# Insert another new empty paragraph after the Biography paragraph and
# turn it into a hyperlink pointing at the bio-examples article.
$bioPara.Range.InsertParagraphAfter()
$linkPara = $bioPara.Next()

$url = "https://blog.hubspot.com/marketing/professional-bio-examples#short"

# Put the display text in first (as a clean single run) ...
$linkPara.Range.Text = $url

# ... then clear it and add the hyperlink at the same (now empty)
# spot, explicitly supplying the display text. This avoids leaving a
# stray empty run behind in the paragraph (which happens if Hyperlinks.Add
# is targeted directly at the freshly-inserted empty paragraph).
$textRange = $d.Range($linkPara.Range.Start, $linkPara.Range.Start + $url.Length)
$textRange.Text = ""

$insertionPoint = $d.Range($linkPara.Range.Start, $linkPara.Range.Start)
$d.Hyperlinks.Add($insertionPoint, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url) | Out-Null
